$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header row (row 1) labels to the new persistent-URI vocabulary
# naming scheme, leaving the data rows (M/F gender codes + translations)
# untouched.
$ws.Range("A1").Value = "codice_1 _Llivello"
$ws.Range("B1").Value = "label_ITA_1_livello"
$ws.Range("C1").Value = "label_ENG_1_livello"
$ws.Range("D1").Value = "label_DEU_1_livello"
$ws.Range("E1").Value = "label_FRA_1_livello"

# Update the saved selection / active cell to match the source workbook
$ws.Range("E2").Select()
